$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo in BOM: row 9 (J4 / OLED display module) had "-" (not applicable)
# in the MOUSER (H) and DIGIKEY (I) columns, but should read "N.M." like the
# LCSC column (J) already does.
$ws.Range("H9").Value = "N.M."
$ws.Range("I9").Value = "N.M."

# Reflect the active cell/selection left behind after the edit.
$ws.Range("I9").Select()
